$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record_PlayerHero")
$ws.Activate()

# New skill columns appended after Talent8 (column AF), i.e. AG..AP.
$ws.Cells.Item(1, 33).Value = "Skill1"
$ws.Cells.Item(1, 34).Value = "Skill2"
$ws.Cells.Item(1, 35).Value = "Skill3"
$ws.Cells.Item(1, 36).Value = "Skill4"
$ws.Cells.Item(1, 37).Value = "Skill5"
$ws.Cells.Item(1, 38).Value = "SkillLevel1"
$ws.Cells.Item(1, 39).Value = "SkillLevel2"
$ws.Cells.Item(1, 40).Value = "SkillLevel3"
$ws.Cells.Item(1, 41).Value = "SkillLevel4"
$ws.Cells.Item(1, 42).Value = "SkillLevel5"

$ws.Cells.Item(2, 33).Value = "string"
$ws.Cells.Item(2, 34).Value = "string"
$ws.Cells.Item(2, 35).Value = "string"
$ws.Cells.Item(2, 36).Value = "string"
$ws.Cells.Item(2, 37).Value = "string"
$ws.Cells.Item(2, 38).Value = "int"
$ws.Cells.Item(2, 39).Value = "int"
$ws.Cells.Item(2, 40).Value = "int"
$ws.Cells.Item(2, 41).Value = "int"
$ws.Cells.Item(2, 42).Value = "int"

# Give the new header cells (AG1:AP1) the same style used by the other
# headers in that row (W1:AF1, style index 15 / applyFont only).
$ws.Range("AG1:AP1").Style = $ws.Range("W1").Style

# Column AL ("SkillLevel1") needs a wider, best-fit column width.
$ws.Columns("AL").ColumnWidth = 10.857142857142858

# SaveInterval sample value bumped from 21 to 31.
$ws.Cells.Item(2, 3).Value = 31

# Restore the view: selection back at B2, no frozen/scrolled left column.
$ws.Range("B2").Select()
